# Remove the "Screenshots or Demo" section (its Heading1 title paragraph
# plus the single placeholder paragraph beneath it) from the document,
# exactly as described by the target diff.

$d = $word.ActiveDocument

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $text = $text.TrimEnd([char]13, [char]7)
    if ($text -eq "Screenshots or Demo") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $headingPara = $d.Paragraphs.Item($targetIndex)
    $bodyPara = $d.Paragraphs.Item($targetIndex + 1)

    $rangeStart = $headingPara.Range.Start
    $rangeEnd = $bodyPara.Range.End

    $d.Range($rangeStart, $rangeEnd).Delete()
}
